$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.111.28"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").Value = "'2.620.62"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'604.16"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "'146.52"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").Value = "'2.619.22"
$ws.Range("E9").Value = "  -0.47%  "

$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").Value = "'5.61"
$ws.Range("E11").Value = "  -1.05%  "

$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").Value = "'0.361"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").Value = "'27.20"
$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("D15").Value = "'3.086.66"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").Value = "'62.955.72"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").Value = "'2.605.01"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "'11.24"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "'4.44"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("D21").Value = "'339.39"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").Value = "'6.81"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").Value = "'66.67"
$ws.Range("E24").Value = "  -3.31%  "

$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'1.62"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").Value = "'1.55"
$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("D27").Value = "'8.68"
$ws.Range("E27").Value = "  +3.38%  "

$ws.Range("E28").Value = "  -1.79%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'7.90"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").Value = "'528.67"
$ws.Range("E31").Value = "  +5.05%  "

$ws.Range("D32").Value = "'2.02"
$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("D34").Value = "'0.0₃0802"
$ws.Range("E34").Value = "  -1.49%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  +14.76%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'170.47"
$ws.Range("E36").Value = "  -2.62%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  -0.43%  "

$ws.Range("E40").Value = "  +6.88%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'170.21"
$ws.Range("E42").Value = "  +2.66%  "

$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "'22.27"
$ws.Range("E44").Value = "  +2.78%  "

$ws.Range("D45").Value = "'0.0568"
$ws.Range("E45").Value = "  +3.94%  "

$ws.Range("D46").Value = "'0.623"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").Value = "'0.0960"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("E48").Value = "  -0.33%  "

$ws.Range("D49").Value = "'18.45"
$ws.Range("E49").Value = "  -1.16%  "

$ws.Range("D50").Value = "'1.77"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("D51").Value = "'11.23"
$ws.Range("E51").Value = "  -1.14%  "
